$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.04 = 28142.45 pesos`n✅ 28142.45 pesos = 7.0 = 959.56 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the Binance/transfi rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 141.999
$wsTasas.Range("O10").Value = 3996.2
$wsTasas.Range("N12").Value = 4018
$wsTasas.Range("O12").Value = 137
